$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16 data (appended entry), mirroring the text-cell formatting of
# the existing rows (col A is blank, all other cells are plain text).
$ws.Range("A16:H16").NumberFormat = "@"

$ws.Range("A16").Value = ""
$ws.Range("B16").Value = "أحمد شريم"
$ws.Range("C16").Value = "8"
$ws.Range("D16").Value = "الصمود"
$ws.Range("E16").Value = "الرحلة 2"
$ws.Range("F16").Value = "C3"
$ws.Range("G16").Value = "IDRF"
$ws.Range("H16").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:١٤:٢٠ م"
